$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column O currently holds "Color 7" (header only, no data beneath it).
# Remove that entire column, shifting "Imágenes" (previously column P) left into O.
$ws.Range("O1").EntireColumn.Delete()
